# Add two more rows of data to the "February" sheet, then move the
# active selection cursor, matching the updated test fixture.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("February")

$ws.Range("A3").Value = 3
$ws.Range("B3").Value = 2021
$ws.Range("A4").Value = 4
$ws.Range("B4").Value = 2022

$ws.Range("C6").Select()
